$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update title/timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 15:15"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5569520
$ws.Range("C4").Value = 2888
$ws.Range("D4").Value = 2922936
$ws.Range("E4").Value = 2473441
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 173143

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2661396
$ws.Range("C6").Value = 14080
$ws.Range("D6").Value = 1932839
$ws.Range("E6").Value = 677381
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 131
$ws.Range("H6").Value = 51176

# Row 16: Arabia Saudita
$ws.Range("A16").Value = "Arabia Saudita"
$ws.Range("B16").Value = 299914
$ws.Range("C16").Value = 1372
$ws.Range("D16").Value = 268385
$ws.Range("E16").Value = 28093
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 3436

# Row 36: Suecia
$ws.Range("A36").Value = "Suecia"
$ws.Range("B36").Value = 85045
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 8
$ws.Range("H36").Value = 5787

# Row 37: China
$ws.Range("A37").Value = "China"
$ws.Range("B37").Value = 84849
$ws.Range("C37").Value = 22
$ws.Range("D37").Value = 79603
$ws.Range("E37").Value = 612
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 4634

# Row 41: Kuwait
$ws.Range("A41").Value = "Kuwait"
$ws.Range("B41").Value = 76827
$ws.Range("C41").Value = 622
$ws.Range("D41").Value = 68633
$ws.Range("E41").Value = 7692
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 502

# Row 44: Emiratos Arabes Unidos
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("B44").Value = 64541
$ws.Range("C44").Value = 229
$ws.Range("D44").Value = 57794
$ws.Range("E44").Value = 6383
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 364

# Row 45: Paises Bajos
$ws.Range("A45").Value = "Paises Bajos"
$ws.Range("B45").Value = 63484
$ws.Range("C45").Value = 482
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 6172

# Row 50: Portugal
$ws.Range("A50").Value = "Portugal"
$ws.Range("B50").Value = 54234
$ws.Range("C50").Value = 132
$ws.Range("D50").Value = 39800
$ws.Range("E50").Value = 12655
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 1779

# Row 60: Afganistan
$ws.Range("A60").Value = "Afganistan"
$ws.Range("B60").Value = 37599
$ws.Range("C60").Value = 3
$ws.Range("D60").Value = 27166
$ws.Range("E60").Value = 9058
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 1375

# Row 64: Kenia
$ws.Range("A64").Value = "Kenia"
$ws.Range("B64").Value = 30365
$ws.Range("C64").Value = 245
$ws.Range("D64").Value = 17160
$ws.Range("E64").Value = 12723
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 8
$ws.Range("H64").Value = 482

# Row 65: Moldavia
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 30183
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 21220
$ws.Range("E65").Value = 8067
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 896

# Row 79: Dinamarca
$ws.Range("A79").Value = "Dinamarca"
$ws.Range("B79").Value = 15740
$ws.Range("C79").Value = 123
$ws.Range("D79").Value = 13417
$ws.Range("E79").Value = 1702
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 621

# Row 98: Finlandia
$ws.Range("A98").Value = "Finlandia"
$ws.Range("B98").Value = 7752
$ws.Range("C98").Value = 21
$ws.Range("D98").Value = 7050
$ws.Range("E98").Value = 368
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 334

# Row 124: Sri Lanka
$ws.Range("A124").Value = "Sri Lanka"
$ws.Range("B124").Value = 2895
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 2676
$ws.Range("E124").Value = 208
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 11

# Row 178: Papua Nueva Guinea
$ws.Range("A178").Value = "Papua Nueva Guinea"
$ws.Range("B178").Value = 333
$ws.Range("C178").Value = 10
$ws.Range("D178").Value = 110
$ws.Range("E178").Value = 220
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 3

# Row 213: Montserrat
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214: Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
